# Case Study - Sales Analysis 2019: remove the "% of total sales" helper
# column and the scratch total/percentage rows from the "q5" sheet, and
# restore the plain selection state on the other lightweight query sheets.

$wb = $excel.ActiveWorkbook

# --- q5: drop the helper column C ("% of total sales") -------------------
# Column C held a formula (=B2/$B$21 etc.) plus its own header/style; a
# whole-column delete shifts D (Price) and E (Qty Ordered) left and fixes
# up the remaining SUM()/ratio formulas' column refs automatically.
$ws5 = $wb.Worksheets.Item("q5")
$ws5.Columns.Item(3).Delete()

# The scratch rows at the bottom (grand total, a spot-check subtotal, and
# a one-off "=B2" echo) go back to being empty - just the label clears and
# the formulas clear, but the pre-existing number formatting (style 4)
# stays on the B cells.
$ws5.Range("A21").ClearContents()
$ws5.Range("B21").ClearContents()
$ws5.Range("B22").ClearContents()
$ws5.Range("B24").ClearContents()

# --- q2 / q3: just a cursor/selection reset, no data changes -------------
$ws2 = $wb.Worksheets.Item("q2")
$ws2.Range("A6").Select()

$ws3 = $wb.Worksheets.Item("q3")
$ws3.Range("A2").Select()

# --- restore q5 as the active tab with its own new selection -------------
$ws5.Activate()
$ws5.Range("B10").Select()
